# Update automàtic: dades i banners [2026-02-20 20:50]
# Applies the scraped meteo.cat data refresh to Dades_Meteo sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "2026-02-20 20:48:29"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "62%"
$ws.Range("N2").Value = "-2.3 °C 20:10 TU"
$ws.Range("O2").Value = "0.5 °C"
# Row 3
$ws.Range("E3").Value = "2026-02-20 20:48:32"
# Row 4
$ws.Range("E4").Value = "2026-02-20 20:48:34"
$ws.Range("J4").Value = "1022.6 hPa"
$ws.Range("O4").Value = "10.2 °C"
# Row 5
$ws.Range("E5").Value = "2026-02-20 20:48:37"
$ws.Range("L5").Value = "23.4 km/h - 342º 20:21 TU"
# Row 6
$ws.Range("E6").Value = "2026-02-20 20:48:39"
$ws.Range("J6").Value = "1022.6 hPa"
# Row 7
$ws.Range("E7").Value = "2026-02-20 20:48:42"
$ws.Range("J7").Value = "1022.5 hPa"
# Row 8
$ws.Range("E8").Value = "2026-02-20 20:48:44"
$ws.Range("J8").Value = "1022.8 hPa"
$ws.Range("O8").Value = "9.3 °C"
# Row 9
$ws.Range("E9").Value = "2026-02-20 20:48:47"
# Row 10
$ws.Range("E10").Value = "2026-02-20 20:48:50"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "79%"
$ws.Range("O10").Value = "7.7 °C"
# Row 11
$ws.Range("E11").Value = "2026-02-20 20:48:52"
# Row 12
$ws.Range("E12").Value = "2026-02-20 20:48:54"
# Row 13
$ws.Range("E13").Value = "2026-02-20 20:48:57"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "44%"
$ws.Range("J13").Value = "1023.6 hPa"
$ws.Range("O13").Value = "6.4 °C"
# Row 14
$ws.Range("E14").Value = "2026-02-20 20:48:59"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "57%"
$ws.Range("N14").Value = "6.4 °C 20:29 TU"
$ws.Range("O14").Value = "12.0 °C"
# Row 15
$ws.Range("E15").Value = "2026-02-20 20:49:02"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "46%"
$ws.Range("O15").Value = "13.3 °C"
# Row 16
$ws.Range("E16").Value = "2026-02-20 20:49:04"
# Row 17
$ws.Range("E17").Value = "2026-02-20 20:49:07"
$ws.Range("M17").Value = "6.0 °C 20:26 TU"
# Row 18
$ws.Range("E18").Value = "2026-02-20 20:49:09"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "77%"
$ws.Range("J18").Value = "1022.9 hPa"
# Row 19
$ws.Range("E19").Value = "2026-02-20 20:49:12"
# Row 20
$ws.Range("E20").Value = "2026-02-20 20:49:14"
$ws.Range("O20").Value = "-2.8 °C"
# Row 21
$ws.Range("E21").Value = "2026-02-20 20:49:17"
$ws.Range("J21").Value = "1022.6 hPa"
# Row 22
$ws.Range("E22").Value = "2026-02-20 20:49:19"
# Row 23
$ws.Range("E23").Value = "2026-02-20 20:49:22"
$ws.Range("O23").Value = "-4.8 °C"
# Row 24
$ws.Range("E24").Value = "2026-02-20 20:49:24"
$ws.Range("J24").Value = "1025.4 hPa"
# Row 25
$ws.Range("E25").Value = "2026-02-20 20:49:27"
$ws.Range("O25").Value = "-1.5 °C"
# Row 26
$ws.Range("E26").Value = "2026-02-20 20:49:29"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "35%"
$ws.Range("J26").Value = "1021.8 hPa"
# Row 27
$ws.Range("E27").Value = "2026-02-20 20:49:32"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "43%"
# Row 28
$ws.Range("E28").Value = "2026-02-20 20:49:34"
$ws.Range("J28").Value = "1023.0 hPa"
$ws.Range("O28").Value = "7.1 °C"
# Row 29
$ws.Range("E29").Value = "2026-02-20 20:49:37"
# Row 30
$ws.Range("E30").Value = "2026-02-20 20:49:39"
$ws.Range("J30").Value = "1022.3 hPa"
# Row 31
$ws.Range("E31").Value = "2026-02-20 20:49:42"
$ws.Range("J31").Value = "1021.5 hPa"
# Row 32
$ws.Range("E32").Value = "2026-02-20 20:49:45"
$ws.Range("O32").Value = "4.5 °C"
# Row 33
$ws.Range("E33").Value = "2026-02-20 20:49:47"
$ws.Range("J33").Value = "1023.1 hPa"
$ws.Range("O33").Value = "6.0 °C"
# Row 34
$ws.Range("E34").Value = "2026-02-20 20:49:50"
$ws.Range("O34").Value = "0.9 °C"
# Row 35
$ws.Range("E35").Value = "2026-02-20 20:49:52"
$ws.Range("J35").Value = "1026.8 hPa"
# Row 36
$ws.Range("E36").Value = "2026-02-20 20:49:55"
$ws.Range("J36").Value = "1022.5 hPa"
# Row 37
$ws.Range("E37").Value = "2026-02-20 20:49:57"
$ws.Range("J37").Value = "1024.5 hPa"
$ws.Range("O37").Value = "4.8 °C"
# Row 38
$ws.Range("E38").Value = "2026-02-20 20:50:00"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "68%"
# Row 39
$ws.Range("E39").Value = "2026-02-20 20:50:02"
# Row 40
$ws.Range("E40").Value = "2026-02-20 20:50:05"
$ws.Range("J40").Value = "1023.4 hPa"
# Row 41
$ws.Range("E41").Value = "2026-02-20 20:50:07"
$ws.Range("J41").Value = "1023.1 hPa"
# Row 42
$ws.Range("E42").Value = "2026-02-20 20:50:10"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "65%"
# Row 43
$ws.Range("E43").Value = "2026-02-20 20:50:12"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "76%"
$ws.Range("O43").Value = "5.0 °C"
# Row 44
$ws.Range("E44").Value = "2026-02-20 20:50:15"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "76%"
# Row 45
$ws.Range("E45").Value = "2026-02-20 20:50:17"
$ws.Range("J45").Value = "1029.6 hPa"
# Row 46
$ws.Range("E46").Value = "2026-02-20 20:50:20"
$ws.Range("J46").Value = "1026.4 hPa"
$ws.Range("N46").Value = "8.5 °C 20:20 TU"
$ws.Range("O46").Value = "12.0 °C"
